# "ControlDbp begin rationalizing commands"
# Swap the "SHOW RUNNING JOB runId" and "STOP JOB runId" command rows on the
# Commands sheet so STOP JOB runId (row 26) now precedes SHOW RUNNING JOB
# runId (row 27). Column D (the URI template) is identical for both rows,
# and column B is a derived formula (LEFT(A,SEARCH(" ",A)-1)), so only the
# command text (A) and HTTP verb (C) need to trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

$commandRow26 = $ws.Range("A26").Value2
$verbRow26    = $ws.Range("C26").Value2

$commandRow27 = $ws.Range("A27").Value2
$verbRow27    = $ws.Range("C27").Value2

$ws.Range("A26").Value = $commandRow27
$ws.Range("C26").Value = $verbRow27

$ws.Range("A27").Value = $commandRow26
$ws.Range("C27").Value = $verbRow26
